# This script applies a within-group reordering of the Nutrient rows for
# the "fish_id_outliers" worksheet (Sheet 1). For a handful of
# Code_sample groups, the rows listing different Nutrient values had been
# entered out of order; this restores the canonical order by rotating the
# Nutrient/concentration/mean_all/min_all/max_all values (columns B-F)
# among the affected rows, while Code_sample (A), Water_percent (G),
# Prepa_operator (H) and Comment (I) stay attached to their original row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=15; B="Ag"; C=0.148680460909429; D=0.053802992568572; E=0.005; F=0.488527763216259}
    @{Row=16; B="Cd"; C=1.94874654111428; D=0.675619363236724; E=0.0482610082961072; F=8.07492237590225}
    @{Row=19; B="Ag"; C=0.235291712662006; D=0.053802992568572; E=0.005; F=0.488527763216259}
    @{Row=20; B="Cd"; C=1.86638865039454; D=0.675619363236724; E=0.0482610082961072; F=8.07492237590225}
    @{Row=23; B="Ag"; C=0.156323737612081; D=0.053802992568572; E=0.005; F=0.488527763216259}
    @{Row=24; B="Cd"; C=2.56606890042473; D=0.675619363236724; E=0.0482610082961072; F=8.07492237590225}
    @{Row=27; B="Ag"; C=0.266022755177212; D=0.053802992568572; E=0.005; F=0.488527763216259}
    @{Row=28; B="Cd"; C=2.11181141032987; D=0.675619363236724; E=0.0482610082961072; F=8.07492237590225}
    @{Row=33; B="Ag"; C=0.30019725486973; D=0.053802992568572; E=0.005; F=0.488527763216259}
    @{Row=34; B="Cd"; C=1.67050217802252; D=0.675619363236724; E=0.0482610082961072; F=8.07492237590225}
    @{Row=39; B="Ag"; C=0.24039613449405; D=0.053802992568572; E=0.005; F=0.488527763216259}
    @{Row=40; B="Cd"; C=1.5294305566648; D=0.675619363236724; E=0.0482610082961072; F=8.07492237590225}
    @{Row=44; B="Ag"; C=0.282718699685655; D=0.053802992568572; E=0.005; F=0.488527763216259}
    @{Row=45; B="Cd"; C=2.53392624396228; D=0.675619363236724; E=0.0482610082961072; F=8.07492237590225}
    @{Row=174; B="Mn"; C=5.60447071121296; D=2.93832418726495; E=0.505013069090625; F=82.9367697861274}
    @{Row=175; B="As"; C=11.1586861255176; D=4.45532071226407; E=0.612445844820796; F=29.8980160943351}
    @{Row=176; B="Sr"; C=368.954287781932; D=79.272744650811; E=13.8326977284032; F=368.954287781932}
    @{Row=212; B="Mn"; C=8.7773100821808; D=2.93832418726495; E=0.505013069090625; F=82.9367697861274}
    @{Row=213; B="As"; C=12.6798957706955; D=4.45532071226407; E=0.612445844820796; F=29.8980160943351}
    @{Row=214; B="Sr"; C=195.241531369012; D=79.272744650811; E=13.8326977284032; F=368.954287781932}
    @{Row=215; B="Mn"; C=12.2320200631017; D=2.93832418726495; E=0.505013069090625; F=82.9367697861274}
    @{Row=216; B="As"; C=17.8060027505865; D=4.45532071226407; E=0.612445844820796; F=29.8980160943351}
    @{Row=217; B="Sr"; C=163.087533371086; D=79.272744650811; E=13.8326977284032; F=368.954287781932}
    @{Row=220; B="Mn"; C=9.44194603434178; D=2.93832418726495; E=0.505013069090625; F=82.9367697861274}
    @{Row=221; B="As"; C=9.46034341782502; D=4.45532071226407; E=0.612445844820796; F=29.8980160943351}
    @{Row=222; B="Sr"; C=199.237530662306; D=79.272744650811; E=13.8326977284032; F=368.954287781932}
    @{Row=239; B="Co"; C=0.174425298384311; D=0.0946495589522645; E=0.0223495702005731; F=0.48519837857757}
    @{Row=240; B="Cd"; C=2.10883831599954; D=0.675619363236724; E=0.0482610082961072; F=8.07492237590225}
    @{Row=242; B="Zn"; C=117.931207923857; D=53.6194671067111; E=17.2797769119108; F=161.07593729738}
    @{Row=243; B="Se"; C=5.85970749825892; D=2.73299221888054; E=1.25943300020396; F=5.85970749825892}
    @{Row=244; B="Sr"; C=185.264644432407; D=79.272744650811; E=13.8326977284032; F=368.954287781932}
    @{Row=263; B="Co"; C=0.246111439259697; D=0.0946495589522645; E=0.0223495702005731; F=0.48519837857757}
    @{Row=264; B="Sr"; C=190.167355778697; D=79.272744650811; E=13.8326977284032; F=368.954287781932}
    @{Row=272; B="Ag"; C=0.488527763216259; D=0.053802992568572; E=0.005; F=0.488527763216259}
    @{Row=273; B="Cd"; C=8.07492237590225; D=0.675619363236724; E=0.0482610082961072; F=8.07492237590225}
)

foreach ($u in $updates) {
    $row = $u.Row
    $ws.Cells.Item($row, 2).Value = $u.B
    $ws.Cells.Item($row, 3).Value = $u.C
    $ws.Cells.Item($row, 4).Value = $u.D
    $ws.Cells.Item($row, 5).Value = $u.E
    $ws.Cells.Item($row, 6).Value = $u.F
}

Write-Host "Applied" $updates.Count "row updates"
